$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45182 (2023-09-13) to 45184 (2023-09-15) for every data row
# (rows 2 through 155).
$ws.Range("C2:C155").Value = 45184
